# Update the Summary sheet with refreshed comparison numbers.
# Column layout (row 1 headers, unchanged):
#   A Organ
#   B Az_Asctb_perfect_matches
#   C Az_unique_CT
#   D ASCTB_unique_CT
#   E Az_cts_not_matched
#   F Az_percent_not_matching(%)
#   G Asctb_cts_not_matched
#   H Asctb_percent_not_matching(%)
#   I Az_missing_cts
#   J Asctb_missing_cts
#   K Az_incorrect_cts
#   L Asctb_incorrect_cts
#   M Az_ct_match_found_crosswalk
#   N Asctb_ct_match_found_crosswalk

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - lung
$ws.Cells.Item(2, 2).Value = 16
$ws.Cells.Item(2, 5).Value = 45
$ws.Cells.Item(2, 6).Value = 104.6511627906977
$ws.Cells.Item(2, 7).Value = 24
$ws.Cells.Item(2, 8).Value = 37.5
$ws.Cells.Item(2, 12).Value = 34
$ws.Cells.Item(2, 13).Value = 6
$ws.Cells.Item(2, 14).Value = 30

# Row 3 - pancreas
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(3, 12).Value = 6

# Row 4 - kidney
$ws.Cells.Item(4, 2).Value = 36
$ws.Cells.Item(4, 5).Value = 17
$ws.Cells.Item(4, 6).Value = 40.47619047619047
$ws.Cells.Item(4, 7).Value = 11
$ws.Cells.Item(4, 8).Value = 23.40425531914894
$ws.Cells.Item(4, 12).Value = 15
$ws.Cells.Item(4, 13).Value = 3

# Row 5 - brain
$ws.Cells.Item(5, 5).Value = 189
$ws.Cells.Item(5, 6).Value = 1890
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 127

# Row 6 - bone_marrow
$ws.Cells.Item(6, 2).Value = 32
$ws.Cells.Item(6, 7).Value = 13
$ws.Cells.Item(6, 8).Value = 35.13513513513514
$ws.Cells.Item(6, 11).Value = 8
$ws.Cells.Item(6, 12).Value = 11

# Row 7 - blood_pmbc
$ws.Cells.Item(7, 2).Value = 29
$ws.Cells.Item(7, 5).Value = 12
$ws.Cells.Item(7, 6).Value = 33.33333333333334
$ws.Cells.Item(7, 7).Value = 14
$ws.Cells.Item(7, 8).Value = 53.84615384615385
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 4
$ws.Cells.Item(7, 13).Value = 39
